$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.457.78"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.825.28"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "'332.47"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.4570"
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("D8").Value = "'0.3798"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").Value = "'46.36"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'0.07866"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "'0.9693"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").Value = "'21.00"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.824.79"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.879"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "'7.038"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "'89.72"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "'0.06646"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "'17.07"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'1.006"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "27.441.08"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'5.331"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").Value = "'10.80"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'2.309"
$ws.Range("D26").Value = "2.038.84"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "'155.67"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'19.37"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").Value = "'2.052"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "'5.272"
$ws.Range("D31").Value = "'118.27"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'0.9414"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").Value = "'0.09302"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "'5.243"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "'1.317"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'0.05920"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("D38").Value = "'0.02185"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").Value = "'8.046"
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("D40").Value = "'1.155"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Value = "'0.5752"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Value = "'0.1824"
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("D43").Value = "'9.962"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "'1.266"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "'12.01"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").Value = "'0.5430"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("D47").Value = "'1.862"
$ws.Range("D48").Value = "'110.73"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").Value = "'0.06605"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "'1.006"
$ws.Range("D51").Value = "'1.041"
$ws.Range("E51").Value = "  -1.24%  "
